$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 21: "headline+snippet整段文字，abstract整段文字的emotion"
$ws.Range("B21").Value = 0.94191199999999997
$ws.Range("C21").Value = 0.92674999999999996
$ws.Range("D21").Value = 0.89998
$ws.Range("C21:D21").Style = "常规"
$ws.Range("C21:D21").VerticalAlignment = -4108
$ws.Range("F21").Value = "headline+snippet整段文字，abstract整段文字的emotion"

# New row 22: "headline整段文字，abstract整段文字的emotion"
$ws.Range("B22").Value = 0.93780249999999998
$ws.Range("C22").Value = 0.92898999999999998
$ws.Range("D22").Value = 0.89517000000000002
$ws.Range("C22:D22").Style = "常规"
$ws.Range("C22:D22").VerticalAlignment = -4108
$ws.Range("F22").Value = "headline整段文字，abstract整段文字的emotion"

# Update the active selection to match the author's final cursor position
$ws.Range("C21:D22").Select() | Out-Null
